$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price observation was recorded for the week, which gets inserted
# as row 430 (pushing the previous rows 430-524 down to 431-525).
$ws.Range("A430:R430").Insert("xlShiftDown")

$ws.Range("A430").Value = 4
$ws.Range("B430").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C430").Value = "Los Lagos"
$ws.Range("D430").Value = 45173
$ws.Range("E430").Value = 10
$ws.Range("F430").Value = 100112040
$ws.Range("G430").Value = "Cilantro"
$ws.Range("H430").Value = "Sin especificar"
$ws.Range("I430").Value = "Primera"
$ws.Range("J430").Value = 70
$ws.Range("K430").Value = 14000
$ws.Range("L430").Value = 14000
$ws.Range("M430").Value = 14000
$ws.Range("N430").Value = "`$/caja 36 atados"
$ws.Range("O430").Value = "Región Metropolitana"
$ws.Range("P430").Value = 389
$ws.Range("Q430").Value = 36
$ws.Range("R430").Value = "Hortaliza"
